$p = $ppt.ActivePresentation
$p.Slides.Item(13).Delete()
